$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = 'sv'
$ws.Range("J11").Value = 'Statement-opinion'
$ws.Range("I54").Value = 'sd'
$ws.Range("J54").Value = 'Statement-non-opinion'
$ws.Range("I70").Value = 'sd'
$ws.Range("J70").Value = 'Statement-non-opinion'
$ws.Range("I77").Value = 'sd'
$ws.Range("J77").Value = 'Statement-non-opinion'
$ws.Range("I79").Value = 'sd'
$ws.Range("J79").Value = 'Statement-non-opinion'
$ws.Range("I81").Value = 'sd'
$ws.Range("J81").Value = 'Statement-non-opinion'
$ws.Range("I98").Value = 'sd'
$ws.Range("J98").Value = 'Statement-non-opinion'
$ws.Range("I109").Value = 'ba'
$ws.Range("J109").Value = 'Appreciation'
$ws.Range("I118").Value = 'sd'
$ws.Range("J118").Value = 'Statement-non-opinion'
$ws.Range("I119").Value = 'sd'
$ws.Range("J119").Value = 'Statement-non-opinion'
$ws.Range("I121").Value = 'aa'
$ws.Range("J121").Value = 'Agree/Accept'
$ws.Range("I138").Value = 'sd'
$ws.Range("J138").Value = 'Statement-non-opinion'
$ws.Range("I143").Value = '%'
$ws.Range("J143").Value = 'Uninterpretable'
$ws.Range("I151").Value = 'sd'
$ws.Range("J151").Value = 'Statement-non-opinion'
$ws.Range("I161").Value = 'sd'
$ws.Range("J161").Value = 'Statement-non-opinion'
$ws.Range("I183").Value = 'sv'
$ws.Range("J183").Value = 'Statement-opinion'
$ws.Range("I185").Value = 'sd'
$ws.Range("J185").Value = 'Statement-non-opinion'
$ws.Range("I186").Value = 'sd'
$ws.Range("J186").Value = 'Statement-non-opinion'
$ws.Range("I187").Value = 'sd'
$ws.Range("J187").Value = 'Statement-non-opinion'
$ws.Range("I192").Value = 'sd'
$ws.Range("J192").Value = 'Statement-non-opinion'
$ws.Range("I195").Value = 'sv'
$ws.Range("J195").Value = 'Statement-opinion'
$ws.Range("I199").Value = 'aa'
$ws.Range("J199").Value = 'Agree/Accept'
$ws.Range("I207").Value = 'ba'
$ws.Range("J207").Value = 'Appreciation'
$ws.Range("I223").Value = 'aa'
$ws.Range("J223").Value = 'Agree/Accept'
$ws.Range("I241").Value = 'sd'
$ws.Range("J241").Value = 'Statement-non-opinion'
$ws.Range("I242").Value = '%'
$ws.Range("J242").Value = 'Uninterpretable'
$ws.Range("I257").Value = 'sv'
$ws.Range("J257").Value = 'Statement-opinion'
$ws.Range("I260").Value = 'aa'
$ws.Range("J260").Value = 'Agree/Accept'
$ws.Range("I278").Value = 'ba'
$ws.Range("J278").Value = 'Appreciation'
$ws.Range("I282").Value = 'b'
$ws.Range("J282").Value = 'Acknowledge (Backchannel)'
$ws.Range("I287").Value = 'sd'
$ws.Range("J287").Value = 'Statement-non-opinion'
$ws.Range("I289").Value = 'ba'
$ws.Range("J289").Value = 'Appreciation'
$ws.Range("I316").Value = 'sv'
$ws.Range("J316").Value = 'Statement-opinion'
$ws.Range("I317").Value = 'sv'
$ws.Range("J317").Value = 'Statement-opinion'
$ws.Range("I321").Value = 'sd'
$ws.Range("J321").Value = 'Statement-non-opinion'
$ws.Range("I328").Value = 'sv'
$ws.Range("J328").Value = 'Statement-opinion'
$ws.Range("I333").Value = 'sv'
$ws.Range("J333").Value = 'Statement-opinion'
$ws.Range("I337").Value = 'sv'
$ws.Range("J337").Value = 'Statement-opinion'
$ws.Range("I368").Value = 'sd'
$ws.Range("J368").Value = 'Statement-non-opinion'
$ws.Range("I379").Value = 'sd'
$ws.Range("J379").Value = 'Statement-non-opinion'
$ws.Range("I386").Value = 'sv'
$ws.Range("J386").Value = 'Statement-opinion'
$ws.Range("I396").Value = 'b'
$ws.Range("J396").Value = 'Acknowledge (Backchannel)'
$ws.Range("I398").Value = 'ba'
$ws.Range("J398").Value = 'Appreciation'
$ws.Range("I409").Value = 'ba'
$ws.Range("J409").Value = 'Appreciation'
$ws.Range("I417").Value = 'sv'
$ws.Range("J417").Value = 'Statement-opinion'
$ws.Range("I425").Value = 'sv'
$ws.Range("J425").Value = 'Statement-opinion'
$ws.Range("I437").Value = 'sv'
$ws.Range("J437").Value = 'Statement-opinion'
$ws.Range("I440").Value = 'sd'
$ws.Range("J440").Value = 'Statement-non-opinion'
$ws.Range("I452").Value = 'sd'
$ws.Range("J452").Value = 'Statement-non-opinion'
$ws.Range("I464").Value = 'aa'
$ws.Range("J464").Value = 'Agree/Accept'
$ws.Range("I468").Value = 'b'
$ws.Range("J468").Value = 'Acknowledge (Backchannel)'
$ws.Range("I471").Value = 'sv'
$ws.Range("J471").Value = 'Statement-opinion'
$ws.Range("I478").Value = 'aa'
$ws.Range("J478").Value = 'Agree/Accept'
$ws.Range("I479").Value = 'ba'
$ws.Range("J479").Value = 'Appreciation'
$ws.Range("I484").Value = 'sd'
$ws.Range("J484").Value = 'Statement-non-opinion'
$ws.Range("I500").Value = 'sd'
$ws.Range("J500").Value = 'Statement-non-opinion'
$ws.Range("I504").Value = 'ba'
$ws.Range("J504").Value = 'Appreciation'
$ws.Range("I519").Value = 'sd'
$ws.Range("J519").Value = 'Statement-non-opinion'
$ws.Range("I521").Value = 'sd'
$ws.Range("J521").Value = 'Statement-non-opinion'
